# "Mean Values and Box Plot Section"
#
# The sheet that held the raw XAI metric experiment results is renamed
# from the default "Sheet1" to "ANCHOR" (it's the ANCHOR-metrics sheet of
# the workbook), and the active selection is moved to J14 - presumably
# where the new "Mean Values and Box Plot" section was about to be built
# out to the right of the existing data table (B:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sheet1" -> "ANCHOR"
$ws.Name = "ANCHOR"

# Make sure we're on that sheet, then move the selection to J14
$ws.Activate()
$ws.Range("J14").Select()
